$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("41:41").Insert()

$ws.Range("A41").Value = 4
$ws.Range("B41").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C41").Value = "Los Lagos"
$ws.Range("D41").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0).AddDays(44529)
$ws.Range("E41").Value = 10
$ws.Range("F41").Value = 100112026
$ws.Range("G41").Value = "Haba"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 80
$ws.Range("K41").Value = 11000
$ws.Range("L41").Value = 11000
$ws.Range("M41").Value = 11000
$ws.Range("N41").Value = '$/saco 25 kilos'
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value = 440
$ws.Range("Q41").Value = 25
$ws.Range("R41").Value = "Hortaliza"
